$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Hjemme passive" values were re-pulled/tweaked; columns B:E on rows 1-3
# are overwritten (row 1 headers become 16/20/16/20, rows 2-3 get the new
# trial data), and the active selection on the sheet shrinks to B1:E3.
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 60.262886232055124
$ws.Range("C2").Value = 60.281544889481289
$ws.Range("D2").Value = 57.114956082693503
$ws.Range("E2").Value = 64.260037904442086

$ws.Range("B3").Value = 42.657777568082231
$ws.Range("C3").Value = 56.825263581964592
$ws.Range("D3").Value = 54.670250529191847
$ws.Range("E3").Value = 82.926792321616205

$ws.Range("B1:E3").Select()
